# Daily attendance processing - 2025-11-27 07:03:55
#
# The "Recorded By" column (G) lists the users/processes that touched a
# session record, separated by ", ". A daily reprocessing pass re-derives
# this list and always places the literal "System" marker(s) at the end of
# the list instead of wherever it happened to land before, while leaving
# the relative order of the other contributors (including any
# case-variant "system" entries) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ur = $ws.UsedRange
$rowCount = $ur.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val.IndexOf("System") -lt 0) { continue }

    $parts = $val -split ", "

    $others = @()
    $systemCount = 0
    foreach ($p in $parts) {
        # Exact, case-sensitive match on the literal "System" token only -
        # "system" (lowercase) is a distinct contributor and keeps its spot.
        if ($p.Equals("System")) {
            $systemCount++
        } else {
            $others += $p
        }
    }

    if ($systemCount -eq 0) { continue }

    $newParts = @()
    foreach ($o in $others) { $newParts += $o }
    for ($i = 0; $i -lt $systemCount; $i++) {
        $newParts += "System"
    }

    $newVal = [string]::Join(", ", $newParts)

    if (-not $newVal.Equals($val)) {
        $cell.Value = $newVal
    }
}
